# edit.ps1
# Re-ran the DecisionTree grid search (adds 4 more trial rows, rows 32-35,
# and a couple of cells in row 4 shift slightly from re-running with the same params),
# then refreshed the DecisionTree row on the 'Average Statistics' sheet so the
# reported averages reflect the freshly re-run numbers.

$wb  = $excel.ActiveWorkbook
$dt  = $wb.Worksheets.Item("DecisionTree")
$avg = $wb.Worksheets.Item("Average Statistics")

# --- Re-ran row 4 -> a few metrics recomputed to slightly different values ---
$dt.Range("J4").Value = 0.3164838316483832
$dt.Range("U4").Value = 0.3380281690140845
$dt.Range("Z4").Value = 0.3164838316483832
$dt.Range("AA4").Value = 0.2829024980409071
$dt.Range("AB4").Value = 0.3349986735807151

# --- Newly-run DecisionTree trials appended as rows 32-35 ---
$row32 = @(50,30,150,13,12,150,41841,13947,0.3545565354556535,0.3145479314547932,0.9579837957983796,1.00516240051624,0.2870004092211158,0.5585626771509554,0.2565821911054137,0.3078111996195458,0.2850029874526986,0.2542869092429946,0.5102040816326531,0.2092492149586069,0.2797909407665505,0.2676737160120846,0.3545565354556535,0.3391974749620864,0.3637693123023797,0.3145479314547932,0.3021593570423439,0.3238684864749717,0.9858776361332687,0.9839728890099751,0.985498762401047,0.9837008942879621,0.1753825463130961,0.2606212446401409,0.2420723124750055,0.2012425656450529,0.1206813309266997)
$row33 = @(50,30,150,12,12,150,41841,13947,0.3574962357496236,0.3198537319853732,0.9412059941205995,0.9911808991180899,0.2948016415868673,0.5678614208053073,0.2582535651494433,0.3453596012342748,0.2165935380933386,0.2408789386401327,0.5398671096345515,0.2081620669406929,0.3218914185639229,0.1817087845968712,0.3574962357496236,0.338070488059514,0.3680698961591991,0.3198537319853732,0.2963929888286317,0.3331131942204167,0.9836815270474328,0.9766532157984606,0.9833947966583783,0.9760284125929439,0.1748231810041467,0.2593638512380261,0.244140727888587,0.2015834827580967,0.1200887571111419)
$row34 = @(50,30,150,15,12,150,41841,13947,0.3545804354580435,0.3166272316627232,0.9425204942520494,0.9819315981931598,0.2897158322056834,0.6029601029601029,0.2370443349753694,0.3696244520791375,0.1224407868325973,0.2427101200686106,0.5575689991636466,0.2003434459072696,0.3387323943661972,0.1109799291617473,0.3545804354580435,0.3197117576777097,0.3728384159466904,0.3166272316627232,0.2830373001711283,0.3351463062286507,0.9797565572870844,0.9789337946707127,0.9793059704233507,0.9794264830728961,0.1733614515419123,0.2574720728701375,0.243424685685736,0.2038567260727902,0.1218850638294205)
$row35 = @(50,30,150,13,12,150,41841,13947,0.3527401352740135,0.3167706316770632,0.9337491933749194,0.9772710977271097,0.246147495872317,0.6185767097966728,0.2307015829587649,0.3865131578947368,0.1250249650489315,0.1951915240423798,0.5654320987654321,0.1935483870967742,0.3658360418923799,0.1228280407429599,0.3527401352740135,0.3154364119471501,0.3733683183442749,0.3167706316770632,0.2814685527386217,0.3375914725609019,0.987461096996439,0.9820436717732028,0.9870073526087542,0.9816712538876716,0.1760814942100116,0.2590932541327399,0.2437585430138087,0.2022762779546864,0.1187904306887337)

$newRows = @($row32, $row33, $row34, $row35)
$startRow = 32
$r = $startRow
foreach ($rowData in $newRows) {
    $arr = New-Object 'object[,]' 1,$rowData.Length
    for ($i = 0; $i -lt $rowData.Length; $i++) { $arr[0,$i] = $rowData[$i] }
    $targetAddr = "A" + $r + ":AK" + $r
    $dt.Range($targetAddr).Value = $arr
    $r = $r + 1
}

# --- Refresh the DecisionTree averages row on 'Average Statistics' ---
$avgRow = @(0.3541924118898294,0.3170911728855879,0.9477271800668358,0.9957275231021638,0.3541924118898295,0.3306877541584304,0.3675521094918608,0.3170911728855879,0.2928692149416022,0.3314318899566348,0.9847581645604119,0.979903450662093,0.9847374514200574,0.979964033103721,0.1742736574254909,0.2591933471236853,0.2442210683618885,0.2023361951488261,0.1199757319401054)
$arr2 = New-Object 'object[,]' 1,$avgRow.Length
for ($i = 0; $i -lt $avgRow.Length; $i++) { $arr2[0,$i] = $avgRow[$i] }
$avg.Range("B2:T2").Value = $arr2

Write-Host "edit.ps1 applied: DecisionTree rows 32-35 added; averages refreshed."
